$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2861.25  # ALC!H40: 3500 -> 2861.25
$ws.Cells.Item(40, 9).Value = 3001  # ALC!I40: 3000.5 -> 3001
$ws.Cells.Item(40, 10).Value = 2814.6667  # ALC!J40: 3999.5 -> 2814.6667
$ws.Cells.Item(40, 11).Value = 3001  # ALC!K40: 3000.5 -> 3001
$ws.Cells.Item(40, 12).Value = 2814.6667  # ALC!L40: 3999.5 -> 2814.6667
$ws.Cells.Item(40, 13).Value = -2826  # ALC!M40: -2825.5 -> -2826
$ws.Cells.Item(40, 14).Value = -3164.6667  # ALC!N40: -4349.5 -> -3164.6667
$ws.Cells.Item(55, 8).Value = 2671  # ALC!H55: 5116.75 -> 2671
$ws.Cells.Item(55, 9).Value = 228.33333  # ALC!I55: 235 -> 228.33333
$ws.Cells.Item(55, 10).Value = 9999  # ALC!J55: 9998.5 -> 9999
$ws.Cells.Item(55, 11).Value = 228.33333  # ALC!K55: 235 -> 228.33333
$ws.Cells.Item(55, 12).Value = 9999  # ALC!L55: 9998.5 -> 9999
$ws.Cells.Item(55, 13).Value = -14.33332999999999  # ALC!M55: -21 -> -14.33332999999999
$ws.Cells.Item(55, 14).Value = -10427  # ALC!N55: -10426.5 -> -10427
$ws.Cells.Item(106, 8).Value = 222225550  # ALC!H106: 142859980 -> 222225550
$ws.Cells.Item(106, 9).Value = 250003180  # ALC!I106: 153848860 -> 250003180
$ws.Cells.Item(106, 11).Value = 250003180  # ALC!K106: 153848860 -> 250003180
$ws.Cells.Item(106, 13).Value = -250002549  # ALC!M106: -153848229 -> -250002549
$ws.Cells.Item(141, 8).Value = 9048.5  # ALC!H141: 7438.8 -> 9048.5
$ws.Cells.Item(141, 9).Value = 3097  # ALC!I141: 2398 -> 3097
$ws.Cells.Item(141, 11).Value = 9291  # ALC!K141: 7194 -> 9291
$ws.Cells.Item(141, 13).Value = -4111  # ALC!M141: -2014 -> -4111
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1443.75  # ARM!H2: 1381 -> 1443.75
$ws.Cells.Item(2, 9).Value = 1443.75  # ARM!I2: 1381 -> 1443.75
$ws.Cells.Item(2, 11).Value = 1443.75  # ARM!K2: 1381 -> 1443.75
$ws.Cells.Item(2, 13).Value = -1330.75  # ARM!M2: -1268 -> -1330.75
$ws.Cells.Item(32, 8).Value = 158658.06  # ARM!H32: 151606.56 -> 158658.06
$ws.Cells.Item(32, 9).Value = 182665.58  # ARM!I32: 173269.5 -> 182665.58
$ws.Cells.Item(32, 10).Value = 11945.444  # ARM!J32: 12001 -> 11945.444
$ws.Cells.Item(32, 11).Value = 182665.58  # ARM!K32: 173269.5 -> 182665.58
$ws.Cells.Item(32, 12).Value = 11945.444  # ARM!L32: 12001 -> 11945.444
$ws.Cells.Item(32, 13).Value = -182378.58  # ARM!M32: -172982.5 -> -182378.58
$ws.Cells.Item(32, 14).Value = -12519.444  # ARM!N32: -12575 -> -12519.444
$ws.Cells.Item(34, 8).Value = 289666.66  # ARM!H34: 254750 -> 289666.66
$ws.Cells.Item(34, 9).Value = 184500  # ARM!I34: 173000 -> 184500
$ws.Cells.Item(34, 11).Value = 184500  # ARM!K34: 173000 -> 184500
$ws.Cells.Item(34, 13).Value = -184229  # ARM!M34: -172729 -> -184229
$ws.Cells.Item(43, 8).Value = 25358  # ARM!H43: 24438.25 -> 25358
$ws.Cells.Item(43, 9).Value = 0  # ARM!I43: 18000 -> 0
$ws.Cells.Item(43, 11).Value = 0  # ARM!K43: 18000 -> 0
$ws.Cells.Item(43, 13).ClearContents()  # ARM!M43: -17687 -> (removed)
$ws.Cells.Item(116, 8).Value = 1443.75  # ARM!H116: 1381 -> 1443.75
$ws.Cells.Item(116, 9).Value = 1443.75  # ARM!I116: 1381 -> 1443.75
$ws.Cells.Item(116, 11).Value = 1443.75  # ARM!K116: 1381 -> 1443.75
$ws.Cells.Item(116, 13).Value = 850.25  # ARM!M116: 913 -> 850.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1443.75  # BSM!H3: 1381 -> 1443.75
$ws.Cells.Item(3, 9).Value = 1443.75  # BSM!I3: 1381 -> 1443.75
$ws.Cells.Item(3, 11).Value = 1443.75  # BSM!K3: 1381 -> 1443.75
$ws.Cells.Item(3, 13).Value = -1329.75  # BSM!M3: -1267 -> -1329.75
$ws.Cells.Item(12, 8).Value = 3790.4  # BSM!H12: 3776 -> 3790.4
$ws.Cells.Item(12, 9).Value = 3517.3333  # BSM!I12: 3776 -> 3517.3333
$ws.Cells.Item(12, 10).Value = 4200  # BSM!J12: 0 -> 4200
$ws.Cells.Item(12, 11).Value = 3517.3333  # BSM!K12: 3776 -> 3517.3333
$ws.Cells.Item(12, 12).Value = 4200  # BSM!L12: 0 -> 4200
$ws.Cells.Item(12, 13).Value = -3349.3333  # BSM!M12: -3608 -> -3349.3333
$ws.Cells.Item(12, 14).Value = -4536  # BSM!N12: None -> -4536
$ws.Cells.Item(88, 8).Value = 30343  # BSM!H88: 20000 -> 30343
$ws.Cells.Item(88, 10).Value = 30343  # BSM!J88: 20000 -> 30343
$ws.Cells.Item(88, 12).Value = 30343  # BSM!L88: 20000 -> 30343
$ws.Cells.Item(88, 14).Value = -31155  # BSM!N88: -20812 -> -31155
$ws.Cells.Item(91, 8).Value = 30343  # BSM!H91: 20000 -> 30343
$ws.Cells.Item(91, 10).Value = 30343  # BSM!J91: 20000 -> 30343
$ws.Cells.Item(91, 12).Value = 30343  # BSM!L91: 20000 -> 30343
$ws.Cells.Item(91, 14).Value = -33151  # BSM!N91: -22808 -> -33151
$ws.Cells.Item(105, 8).Value = 10446.348  # BSM!H105: 10463.305 -> 10446.348
$ws.Cells.Item(105, 9).Value = 10907.75  # BSM!I105: 11798.546 -> 10907.75
$ws.Cells.Item(105, 10).Value = 9943  # BSM!J105: 9239.333000000001 -> 9943
$ws.Cells.Item(105, 11).Value = 10907.75  # BSM!K105: 11798.546 -> 10907.75
$ws.Cells.Item(105, 12).Value = 9943  # BSM!L105: 9239.333000000001 -> 9943
$ws.Cells.Item(105, 13).Value = -9160.75  # BSM!M105: -10051.546 -> -9160.75
$ws.Cells.Item(105, 14).Value = -13437  # BSM!N105: -12733.333 -> -13437
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 1177448.1  # CRP!H3: 1401437.8 -> 1177448.1
$ws.Cells.Item(3, 9).Value = 1408937.8  # CRP!I3: 1401437.8 -> 1408937.8
$ws.Cells.Item(3, 10).Value = 20000  # CRP!J3: 0 -> 20000
$ws.Cells.Item(3, 11).Value = 1408937.8  # CRP!K3: 1401437.8 -> 1408937.8
$ws.Cells.Item(3, 12).Value = 20000  # CRP!L3: 0 -> 20000
$ws.Cells.Item(3, 13).Value = -1408824.8  # CRP!M3: -1401324.8 -> -1408824.8
$ws.Cells.Item(3, 14).Value = -20226  # CRP!N3: None -> -20226
$ws.Cells.Item(107, 8).Value = 1527.5264  # CRP!H107: 1568.2222 -> 1527.5264
$ws.Cells.Item(107, 9).Value = 1415.0667  # CRP!I107: 1459.3572 -> 1415.0667
$ws.Cells.Item(107, 11).Value = 1415.0667  # CRP!K107: 1459.3572 -> 1415.0667
$ws.Cells.Item(107, 13).Value = 504.9332999999999  # CRP!M107: 460.6428000000001 -> 504.9332999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 23810218  # CUL!H107: 23256514 -> 23810218
$ws.Cells.Item(107, 10).Value = 47619836  # CUL!J107: 45455340 -> 47619836
$ws.Cells.Item(107, 12).Value = 142859508  # CUL!L107: 136366020 -> 142859508
$ws.Cells.Item(107, 14).Value = -142863348  # CUL!N107: -136369860 -> -142863348
$ws.Cells.Item(129, 8).Value = 3754.7896  # CUL!H129: 4244.579 -> 3754.7896
$ws.Cells.Item(129, 10).Value = 8074.143  # CUL!J129: 9403.571 -> 8074.143
$ws.Cells.Item(129, 12).Value = 24222.429  # CUL!L129: 28210.713 -> 24222.429
$ws.Cells.Item(129, 14).Value = -34222.429  # CUL!N129: -38210.713 -> -34222.429
$ws.Cells.Item(140, 8).Value = 2992.5  # CUL!H140: 3109.7058 -> 2992.5
$ws.Cells.Item(140, 9).Value = 2305.5  # CUL!I140: 2424.182 -> 2305.5
$ws.Cells.Item(140, 11).Value = 6916.5  # CUL!K140: 7272.545999999999 -> 6916.5
$ws.Cells.Item(140, 13).Value = -1736.5  # CUL!M140: -2092.545999999999 -> -1736.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 29990  # GSM!H18: 30000 -> 29990
$ws.Cells.Item(18, 9).Value = 29990  # GSM!I18: 30000 -> 29990
$ws.Cells.Item(18, 11).Value = 29990  # GSM!K18: 30000 -> 29990
$ws.Cells.Item(18, 13).Value = -29697  # GSM!M18: -29707 -> -29697
$ws.Cells.Item(42, 8).Value = 71427  # GSM!H42: 71999.664 -> 71427
$ws.Cells.Item(42, 10).Value = 71427  # GSM!J42: 71999.664 -> 71427
$ws.Cells.Item(42, 12).Value = 71427  # GSM!L42: 71999.664 -> 71427
$ws.Cells.Item(42, 14).Value = -72397  # GSM!N42: -72969.664 -> -72397
$ws.Cells.Item(115, 8).Value = 71427  # GSM!H115: 71999.664 -> 71427
$ws.Cells.Item(115, 10).Value = 71427  # GSM!J115: 71999.664 -> 71427
$ws.Cells.Item(115, 12).Value = 71427  # GSM!L115: 71999.664 -> 71427
$ws.Cells.Item(115, 14).Value = -73777  # GSM!N115: -74349.664 -> -73777
$ws.Cells.Item(132, 8).Value = 1075323.6  # GSM!H132: 921969.75 -> 1075323.6
$ws.Cells.Item(132, 9).Value = 9323.267  # GSM!I132: 7749.2104 -> 9323.267
$ws.Cells.Item(132, 11).Value = 27969.801  # GSM!K132: 23247.6312 -> 27969.801
$ws.Cells.Item(132, 13).Value = -25439.801  # GSM!M132: -20717.6312 -> -25439.801
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(26, 8).Value = 0  # LTW!H26: 10000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # LTW!J26: 10000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # LTW!L26: 10000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # LTW!N26: -10590 -> (removed)
$ws.Cells.Item(93, 8).Value = 1222.375  # LTW!H93: 1286.5 -> 1222.375
$ws.Cells.Item(93, 9).Value = 1271.4166  # LTW!I93: 1356.9166 -> 1271.4166
$ws.Cells.Item(93, 11).Value = 1271.4166  # LTW!K93: 1356.9166 -> 1271.4166
$ws.Cells.Item(93, 13).Value = -23.41660000000002  # LTW!M93: -108.9166 -> -23.41660000000002
$ws.Cells.Item(121, 8).Value = 97998.664  # LTW!H121: 97999 -> 97998.664
$ws.Cells.Item(121, 10).Value = 97998.664  # LTW!J121: 97999 -> 97998.664
$ws.Cells.Item(121, 12).Value = 97998.664  # LTW!L121: 97999 -> 97998.664
$ws.Cells.Item(121, 14).Value = -101492.664  # LTW!N121: -101493 -> -101492.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 59999  # WVR!H54: 0 -> 59999
$ws.Cells.Item(54, 10).Value = 59999  # WVR!J54: 0 -> 59999
$ws.Cells.Item(54, 12).Value = 59999  # WVR!L54: 0 -> 59999
$ws.Cells.Item(54, 14).Value = -61039  # WVR!N54: None -> -61039
$ws.Cells.Item(69, 8).Value = 50901.5  # WVR!H69: 49271 -> 50901.5
$ws.Cells.Item(69, 9).Value = 50000  # WVR!I69: 0 -> 50000
$ws.Cells.Item(69, 10).Value = 51202  # WVR!J69: 49271 -> 51202
$ws.Cells.Item(69, 11).Value = 50000  # WVR!K69: 0 -> 50000
$ws.Cells.Item(69, 12).Value = 51202  # WVR!L69: 49271 -> 51202
$ws.Cells.Item(69, 13).Value = -49251  # WVR!M69: None -> -49251
$ws.Cells.Item(69, 14).Value = -52700  # WVR!N69: -50769 -> -52700
$ws.Cells.Item(72, 8).Value = 50901.5  # WVR!H72: 49271 -> 50901.5
$ws.Cells.Item(72, 9).Value = 50000  # WVR!I72: 0 -> 50000
$ws.Cells.Item(72, 10).Value = 51202  # WVR!J72: 49271 -> 51202
$ws.Cells.Item(72, 11).Value = 150000  # WVR!K72: 0 -> 150000
$ws.Cells.Item(72, 12).Value = 153606  # WVR!L72: 147813 -> 153606
$ws.Cells.Item(72, 13).Value = -146256  # WVR!M72: None -> -146256
$ws.Cells.Item(72, 14).Value = -161094  # WVR!N72: -155301 -> -161094
$ws.Cells.Item(96, 8).Value = 2972.0908  # WVR!H96: 2607.1428 -> 2972.0908
$ws.Cells.Item(96, 9).Value = 3024.5  # WVR!I96: 2962.375 -> 3024.5
$ws.Cells.Item(96, 10).Value = 2832.3333  # WVR!J96: 2133.5 -> 2832.3333
$ws.Cells.Item(96, 11).Value = 3024.5  # WVR!K96: 2962.375 -> 3024.5
$ws.Cells.Item(96, 12).Value = 2832.3333  # WVR!L96: 2133.5 -> 2832.3333
$ws.Cells.Item(96, 13).Value = -1651.5  # WVR!M96: -1589.375 -> -1651.5
$ws.Cells.Item(96, 14).Value = -5578.3333  # WVR!N96: -4879.5 -> -5578.3333
$ws.Cells.Item(113, 8).Value = 413.7619  # WVR!H113: 420.5238 -> 413.7619
$ws.Cells.Item(113, 10).Value = 535.7143  # WVR!J113: 556 -> 535.7143
$ws.Cells.Item(113, 12).Value = 1607.1429  # WVR!L113: 1668 -> 1607.1429
$ws.Cells.Item(113, 14).Value = -5947.1429  # WVR!N113: -6008 -> -5947.1429
$ws.Cells.Item(129, 8).Value = 37328  # WVR!H129: 0 -> 37328
$ws.Cells.Item(129, 10).Value = 37328  # WVR!J129: 0 -> 37328
$ws.Cells.Item(129, 12).Value = 37328  # WVR!L129: 0 -> 37328
$ws.Cells.Item(129, 14).Value = -47328  # WVR!N129: None -> -47328
$ws.Cells.Item(132, 8).Value = 1968.8628  # WVR!H132: 2013.54 -> 1968.8628
$ws.Cells.Item(132, 9).Value = 1770.0526  # WVR!I132: 1825.0541 -> 1770.0526
$ws.Cells.Item(132, 11).Value = 5310.1578  # WVR!K132: 5475.1623 -> 5310.1578
$ws.Cells.Item(132, 13).Value = -2780.1578  # WVR!M132: -2945.1623 -> -2780.1578
